$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (title) gains a thin box border around A1:I1, font stays the same (bold 16pt 標楷體) ---
$ws.Range("A2").Copy()
$ws.Range("A1:I1").PasteSpecial(-4122)
$titleRng = $ws.Range("A1:I1")
$titleRng.Font.Bold = $true
$titleRng.Font.Size = 16
$titleRng.Font.Name = "標楷體"

# --- Row 3: copy row 2's formatting/height down onto row 3, then type a blank placeholder into A3 ---
$ws.Range("A2:I2").Copy()
$ws.Range("A3:I3").PasteSpecial(-4122)
$ws.Rows(3).RowHeight = $ws.Rows(2).RowHeight
$ws.Range("A3").Value = "  "

# --- Row 4: copy the plain (unformatted) column-A style down across A4:I4, then fill with blank placeholders ---
$ws.Range("A10").Copy()
$ws.Range("A4:I4").PasteSpecial(-4122)
$ws.Range("A4:I4").Value = "  "

# --- Remove the now unused conditional formatting on the 職級 columns ---
$ws.Range("G3:H3").FormatConditions.Delete()
$ws.Range("I3").FormatConditions.Delete()
$ws.Range("G4:H4").FormatConditions.Delete()
$ws.Range("I4").FormatConditions.Delete()

$ws.Application.CutCopyMode = $false

# --- Restore selection to match the saved cursor position ---
$ws.Range("G7").Select() | Out-Null
